$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 158, shifting rows 158:251 down to 159:252
$ws.Rows.Item(158).Insert()

$ws.Cells.Item(158, 1).Value = 10
$ws.Cells.Item(158, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(158, 3).Value = "La Araucanía"
$ws.Cells.Item(158, 4).Value = 45001
$ws.Cells.Item(158, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(158, 5).Value = 9
$ws.Cells.Item(158, 6).Value = 100112012
$ws.Cells.Item(158, 7).Value = "Espinaca"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 100
$ws.Cells.Item(158, 11).Value = 10000
$ws.Cells.Item(158, 12).Value = 10000
$ws.Cells.Item(158, 13).Value = 10000
$ws.Cells.Item(158, 14).Value = "`$/docena de atados"
$ws.Cells.Item(158, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(158, 16).Value = 3333
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
